# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) for the last data row
# (the d80a9b47-... entry, row 5) on both the zh-cn and de-de sheets,
# recording the timestamp at which the handoff report was generated.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-02 05:53:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-02 05:53:58"
